# Added WebView handling capability
# - Adds a new "Login" test scenario row to the "Test Info" sheet, changing
#   the platform from iOS to Android for the existing row and the new one.
# - Leaves "Test Data" sheet content unchanged (only selection moves).

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("Test Info")
$wsData = $wb.Worksheets.Item("Test Data")

# Update existing data row's platform from iOS -> Android
$wsInfo.Range("F2").Value = "Android"

# Add the new "Login" scenario row
$wsInfo.Range("A3").Value = "LoginPageTest"
$wsInfo.Range("B3").Value = "Login"
$wsInfo.Range("C3").Value = "To check login func"
$wsInfo.Range("D3").Value = "Yes"
$wsInfo.Range("E3").Value = "Chrome"
$wsInfo.Range("F3").Value = "Android"

# Update selections to match the authored state
$wsInfo.Activate()
$wsInfo.Range("B3").Select()

$wsData.Activate()
$wsData.Range("K25").Select()

$wsInfo.Activate()
